$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Industries") was 1 for rows 22-80 (dates 3/21/2020-5/18/2020).
# Updated policy data sets these to 0.
$ws.Range("H22:H80").Value = 0
